$wb = $excel.ActiveWorkbook
$chr = $wb.Worksheets.Item("Character")
$stage = $wb.Worksheets.Item("Stage")

# ---------------------------------------------------------------------------
# Sheet "Character": add jump-related fields
#   - remove the old "jumpPower" field (column E)
#   - append "jumpMaxCount" (type "int") right after "weight"
#   - insert "width" / "height" fields before "weight"
#   - insert "jumpHeight" (type "float") before "jumpMaxCount"
# ---------------------------------------------------------------------------

# Drop the old jumpPower column entirely (field name + type + value)
$chr.Range("E1").EntireColumn.Delete() | Out-Null
# weight is now column D

# Append a trailing column for jumpMaxCount (type int)
$chr.Range("D1:D5").Copy() | Out-Null
$chr.Range("E1:E5").PasteSpecial(-4122) | Out-Null
$chr.Range("E4").Value = "int"
$chr.Range("E1").Value = "jumpMaxCount"
$chr.Range("E5").Value = 2

# Insert two new blank columns before weight (currently column D) for width/height;
# weight shifts D->F, jumpMaxCount shifts E->G
$chr.Range("D1:E1").EntireColumn.Insert() | Out-Null
$chr.Range("D1").Value = "width"
$chr.Range("E1").Value = "height"
$chr.Range("D4").Value = "float"
$chr.Range("E4").Value = "float"
$chr.Range("D5").Value = 100
$chr.Range("E5").Value = 150
$chr.Range("F5").Value = 100

# ---------------------------------------------------------------------------
# Sheet "Stage": add a new column C "fieldGravity"
# ---------------------------------------------------------------------------

# C1 header -> copy formatting from B1 (yellow header style) then set text
$stage.Range("B1").Copy() | Out-Null
$stage.Range("C1").PasteSpecial(-4122) | Out-Null

# C2 -> copy formatting from B2 (blank bordered cell)
$stage.Range("B2").Copy() | Out-Null
$stage.Range("C2").PasteSpecial(-4122) | Out-Null

# C3 -> enum description, wrapped text, taller row
$stage.Range("B1").Copy() | Out-Null
$stage.Range("C3").PasteSpecial(-4122) | Out-Null
$stage.Range("C3").WrapText = $true
$stage.Range("C3").Value = "None = 0,`nNormal = 1,"
$stage.Rows.Item(3).RowHeight = 34

$stage.Range("C1").Value = "fieldGravity"

# C4 -> type name for the new field
$stage.Range("B4").Copy() | Out-Null
$stage.Range("C4").PasteSpecial(-4122) | Out-Null
$stage.Range("C4").Value = "FieldGravity"

# C5 -> sample value
$stage.Range("B5").Copy() | Out-Null
$stage.Range("C5").PasteSpecial(-4122) | Out-Null
$stage.Range("C5").Value = 1

# widen the new column
$stage.Columns.Item(3).ColumnWidth = 18.6667

# ---------------------------------------------------------------------------
# Back to "Character": insert jumpHeight (type float) right before jumpMaxCount
# (jumpMaxCount shifts G->H)
# ---------------------------------------------------------------------------
$chr.Range("G1").EntireColumn.Insert() | Out-Null
$chr.Range("F1:F5").Copy() | Out-Null
$chr.Range("G1:G5").PasteSpecial(-4122) | Out-Null
$chr.Range("G1").Value = "jumpHeight"
$chr.Range("G4").Value = "float"
$chr.Range("G5").Value = 200

# match column widths of the new trailing columns to the existing data columns
$chr.Columns.Item(7).ColumnWidth = $chr.Columns.Item(3).ColumnWidth
$chr.Columns.Item(8).ColumnWidth = $chr.Columns.Item(3).ColumnWidth

# ---------------------------------------------------------------------------
# Restore selections similar to what Excel would leave after this edit
# ---------------------------------------------------------------------------
$stage.Range("A5:B5").Select() | Out-Null
$chr.Range("H22").Select() | Out-Null
$chr.Activate() | Out-Null
